$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 74, shifting existing rows 74:160 down to 75:161
$ws.Rows("74:74").Insert()

# Populate the newly-inserted row 74 with its data
$ws.Range("A74").Value = 8
$ws.Range("B74").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C74").Value = 'Coquimbo'
$ws.Range("D74").Value = 44494
$ws.Range("E74").Value = 4
$ws.Range("F74").Value = 100112003
$ws.Range("G74").Value = 'Ajo'
$ws.Range("H74").Value = 'Chino'
$ws.Range("I74").Value = 'Primera'
$ws.Range("J74").Value = 440
$ws.Range("K74").Value = 19000
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = 19500
$ws.Range("N74").Value = '$/caja 10 kilos'
$ws.Range("O74").Value = 'China'
$ws.Range("P74").Value = 1950
$ws.Range("Q74").Value = 10
$ws.Range("R74").Value = 'Hortaliza'

# Make sure the new date cell keeps the same date number format as the rest of column D
$ws.Range("D74").NumberFormat = $ws.Range("D75").NumberFormat
